# Add "2022-Q3" data: new quarter sheet + summary row on "总计".
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet right before the existing "2022-Q2" sheet and
#    name it "2022-Q3". (All sheets after it shift right automatically.)
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# Copy look & feel (column layout / header style) from the 2022-Q2 sheet so
# the new sheet matches the existing per-quarter sheet formatting.
$headerSrc = $q2Sheet.Range("A1:H1")
$headerSrc.Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)  # xlPasteFormats

$dataStyleSrc = $q2Sheet.Range("A2:H2")

# ---------------------------------------------------------------------------
# 2) Fill in the header row + the 2022-Q3 holdings table.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$rows = @(
    @(0, "320003", "诺安先锋混合A",            "40.90", "76.11", "5.25", "2.1472", 3),
    @(1, "519772", "交银新生活力灵活配置混合", "49.60", "83.55", "2.69", "1.3342", 10),
    @(2, "400003", "东方精选混合",              "9.93",  "86.94", "4.86", "0.4826", 6),
    @(3, "400001", "东方龙混合",                "2.59",  "85.52", "4.88", "0.1264", 4),
    @(4, "012621", "诺安先锋混合C",             "1.13",  "76.11", "5.25", "0.0593", 3),
    @(5, "000646", "华润元大量化优选混合A",     "1.47",  "73.62", "3.83", "0.0563", 10),
    @(6, "007827", "华润元大量化优选混合C",     "0.19",  "73.62", "3.83", "0.0073", 10),
    @(7, "005281", "中科沃土转型升级灵活配置混合", "0.10", "57.70", "2.83", "0.0028", 9)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2

    # New cells beyond the copied header/style range need the same per-cell
    # formatting the other quarter sheets use (bold index column, border).
    $dataStyleSrc.Copy()
    $newSheet.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)  # xlPasteFormats

    $newSheet.Range("A" + $r).Value = $row[0]
    $newSheet.Range("B" + $r).Value = $row[1]
    $newSheet.Range("C" + $r).Value = $row[2]
    $newSheet.Range("D" + $r).Value = $row[3]
    $newSheet.Range("E" + $r).Value = $row[4]
    $newSheet.Range("F" + $r).Value = $row[5]
    $newSheet.Range("G" + $r).Value = $row[6]
    $newSheet.Range("H" + $r).Value = $row[7]
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push the existing rows down by one and
#    insert the new 2022-Q3 totals at row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$oldLast = $total.Range("A6:D6")
$oldLast.Copy()
$total.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats, grabs row-6 styling for the new row 7

$summaryRows = @(
    @("2022-Q3", 8, 4.22),
    @("2022-Q2", 4, 2.15),
    @("2022-Q1", 6, 2.64),
    @("2021-Q4", 7, 4.33),
    @("2021-Q3", 1, 1.55),
    @("2020-Q4", 9, 2.44)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $total.Range("A" + $r).Value = $i
    $total.Range("B" + $r).Value = $summaryRows[$i][0]
    $total.Range("C" + $r).Value = $summaryRows[$i][1]
    $total.Range("D" + $r).Value = $summaryRows[$i][2]
}

$total.Range("A1").Select()
